# Generate Report for Handoff
# The c0484e70-5f0d-497b-af67-4e02a7f00f69.md entry moves from
# "Handed back: in sync with en-US" to "Ready for handoff" in both locales,
# and a fresh handoff timestamp is recorded for each locale's XLF file.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the c0484e70... entry -----------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 3 is the c0484e70... entry ---------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-08 02:17:17"

# --- de-de sheet: row 3 is the c0484e70... entry ---------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-08 02:17:25"
